# Update cryptocurrency price/volume data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.624.74"
$ws.Range("E2").Value = "  +2.87%  "
$ws.Range("D3").Value = "1.787.01"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D5").Value = "222.89"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D6").Value = "0.553"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D8").Value = "32.16"
$ws.Range("E8").Value = "  +7.27%  "
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("E10").Value = "  +3.88%  "
$ws.Range("D11").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D11").Value = "0.0934"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").Value = "2.044.95"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").Value = "1.789.49"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D14").Value = "10.85"
$ws.Range("E14").Value = "  +7.85%  "
$ws.Range("D15").Value = "34.616.14"
$ws.Range("E15").Value = "  +2.94%  "
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("E17").Value = "  +2.82%  "
$ws.Range("D18").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D18").Value = "68.22"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D19").Value = "252.17"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("E20").Value = "  +6.66%  "
$ws.Range("D21").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D22").Value = "10.38"
$ws.Range("E22").Value = "  +1.45%  "
$ws.Range("D23").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D23").Value = "4.14"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D25").Value = "158.29"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D26").Value = "16.30"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("D31").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D31").Value = "3.72"
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("D34").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D34").Value = "1.86"
$ws.Range("E34").Value = "  +2.39%  "
$ws.Range("D35").Value = "1.425.02"
$ws.Range("E35").Value = "  -3.66%  "
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D37").Value = "0.0188"
$ws.Range("E37").Value = "  +2.44%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D38").Value = "0.627"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D39").Value = "82.75"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  +4.04%  "
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("D43").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D43").Value = "2.05"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D45").Value = "0.0500"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("D46").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D46").Value = "5.90"
$ws.Range("E46").Value = "  +3.85%  "
$ws.Range("D47").Value = "1.942.70"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D49").Value = "103.35"
$ws.Range("E49").Value = "  +6.64%  "
$ws.Range("D50").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D50").Value = "11.86"
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("D51").NumberFormat = "@"  # keep as text, value would otherwise parse as a number
$ws.Range("D51").Value = "49.43"
$ws.Range("E51").Value = "  -2.88%  "
